$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 ("Marking"): Right mark 5 -> 4, Wrong mark -1 -> -2
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# Row 12 ("Total"): Right total 40 -> 32, Wrong total -5 -> -10
$ws.Range("B12").Value = 32
$ws.Range("C12").Value = -10

# E12 score summary text
$ws.Range("E12").Value = "22 / 112"
